$wb = $excel.ActiveWorkbook

# bosqueu
$ws = $wb.Worksheets.Item("bosqueu")
$ws.Range("A2").Value = 108
$ws.Range("B2").Value = 0.6121212121212121
$ws.Range("C2").Value = 0.5668831168831169
$ws.Range("D2").Value = 0.6121212121212121
$ws.Range("E2").Value = 0.5456423347490885
$ws.Range("F2").Value = 0.7710927456382002

# arbolts
$ws = $wb.Worksheets.Item("arbolts")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 0.8885700370022505
$ws.Range("C2").Value = 1.40324588947056
$ws.Range("D2").Value = 1.184586801154968
$ws.Range("E2").Value = -0.6930513021291713

# bosquets
$ws = $wb.Worksheets.Item("bosquets")
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = 0.6327999337790693
$ws.Range("C2").Value = 0.6461734745561803
$ws.Range("D2").Value = 0.8038491615696196
$ws.Range("E2").Value = 0.2203755231298511

# knnts
$ws = $wb.Worksheets.Item("knnts")
$ws.Range("A2").Value = 23
$ws.Range("B2").Value = 0.6465181857472145
$ws.Range("C2").Value = 0.7076411505920059
$ws.Range("D2").Value = 0.8412140931962599
$ws.Range("E2").Value = 0.1462132328768058

# arboltd
$ws = $wb.Worksheets.Item("arboltd")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 0.6798361056886949
$ws.Range("C2").Value = 0.7604864217832062
$ws.Range("D2").Value = 0.8720587261092032
$ws.Range("E2").Value = 0.3001077355661667

# bosquetd
$ws = $wb.Worksheets.Item("bosquetd")
$ws.Range("A2").Value = 122
$ws.Range("B2").Value = 0.4583290788826664
$ws.Range("C2").Value = 0.410241144989027
$ws.Range("D2").Value = 0.6405006986639648
$ws.Range("E2").Value = 0.6224461138214116

# arbolcc
$ws = $wb.Worksheets.Item("arbolcc")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 0.5664313171976088
$ws.Range("C2").Value = 0.5982123173087937
$ws.Range("D2").Value = 0.7734418642075135
$ws.Range("E2").Value = 0.4124428918192228

# bosquecc
$ws = $wb.Worksheets.Item("bosquecc")
$ws.Range("A2").Value = 115
$ws.Range("B2").Value = 0.4295092022072674
$ws.Range("C2").Value = 0.3504380316051455
$ws.Range("D2").Value = 0.5919780668277715
$ws.Range("E2").Value = 0.6558038835228504

# arbolpp
$ws = $wb.Worksheets.Item("arbolpp")
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = 0.7708994618676629
$ws.Range("C2").Value = 1.875120075338438
$ws.Range("D2").Value = 1.369350238375281
$ws.Range("E2").Value = -1.171224911562941

# bosquepp
$ws = $wb.Worksheets.Item("bosquepp")
$ws.Range("A2").Value = 110
$ws.Range("B2").Value = 0.6279822259756672
$ws.Range("C2").Value = 0.5739241681786854
$ws.Range("D2").Value = 0.7575778297829771
$ws.Range("E2").Value = 0.3354460507961389
